$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")
$ws.Rows.Item(33).Insert()
$ws.Range("R33").Value = "your relationship"
$ws.Range("S33").Value = "2024-09-06 12:23:25"
